$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new weekly price records arrived; insert two blank rows above the
# existing row 277 so the prior rows 277-280 shift down to 279-282 intact,
# then populate the freshly inserted rows 277-278 with the new data.
$ws.Range("A277:R278").EntireRow.Insert()

# New row 277 (Betarraga, Región Metropolitana, $/malla 15 kilos)
$ws.Cells.Item(277, 1).Value = 4
$ws.Cells.Item(277, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(277, 3).Value = "Los Lagos"
$ws.Cells.Item(277, 4).Value = 44656
$ws.Cells.Item(277, 5).Value = 10
$ws.Cells.Item(277, 6).Value = 100114014
$ws.Cells.Item(277, 7).Value = "Betarraga"
$ws.Cells.Item(277, 8).Value = "Sin especificar"
$ws.Cells.Item(277, 9).Value = "Primera"
$ws.Cells.Item(277, 10).Value = 120
$ws.Cells.Item(277, 11).Value = 11000
$ws.Cells.Item(277, 12).Value = 11000
$ws.Cells.Item(277, 13).Value = 11000
$ws.Cells.Item(277, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(277, 15).Value = "Región Metropolitana"
$ws.Cells.Item(277, 16).Value = 733
$ws.Cells.Item(277, 17).Value = 15
$ws.Cells.Item(277, 18).Value = "Hortaliza"

# New row 278 (Betarraga, Región del Maule, $/paquete 5 unidades)
$ws.Cells.Item(278, 1).Value = 4
$ws.Cells.Item(278, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(278, 3).Value = "Los Lagos"
$ws.Cells.Item(278, 4).Value = 44656
$ws.Cells.Item(278, 5).Value = 10
$ws.Cells.Item(278, 6).Value = 100114014
$ws.Cells.Item(278, 7).Value = "Betarraga"
$ws.Cells.Item(278, 8).Value = "Sin especificar"
$ws.Cells.Item(278, 9).Value = "Primera"
$ws.Cells.Item(278, 10).Value = 800
$ws.Cells.Item(278, 11).Value = 1000
$ws.Cells.Item(278, 12).Value = 1000
$ws.Cells.Item(278, 13).Value = 1000
$ws.Cells.Item(278, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(278, 15).Value = "Región del Maule"
$ws.Cells.Item(278, 16).Value = 200
$ws.Cells.Item(278, 17).Value = 5
$ws.Cells.Item(278, 18).Value = "Hortaliza"
